$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 32   Number  21"
$ws.Range("C9").Value = "Report Covering the Week  5/19/2025  Through  5/25/2025"

# --- Crime statistics table updates (rows 15-31) ---
# Row 15
$ws.Range("D15").Value = 2
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("E15").Value = -50
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G15").Value = 2
$ws.Range("G15").NumberFormat = '#,##0'
$ws.Range("H15").Value = 100
$ws.Range("H15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("I15").Value = 18
$ws.Range("J15").Value = 14
$ws.Range("K15").Value = 28.571428571428
$ws.Range("L15").Value = 200
$ws.Range("M15").Value = 800
$ws.Range("N15").Value = 80

# Row 16
$ws.Range("C16").Value = 9
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = 28.571428571428
$ws.Range("F16").Value = 34
$ws.Range("G16").Value = 30
$ws.Range("H16").Value = 13.333333333333
$ws.Range("I16").Value = 140
$ws.Range("J16").Value = 147
$ws.Range("K16").Value = -4.761904761904
$ws.Range("L16").Value = -30.348258706467
$ws.Range("M16").Value = 133.333333333333
$ws.Range("N16").Value = -86.069651741293

# Row 17
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 38
$ws.Range("G17").Value = 45
$ws.Range("H17").Value = -15.555555555555
$ws.Range("I17").Value = 224
$ws.Range("J17").Value = 220
$ws.Range("K17").Value = 1.818181818181
$ws.Range("L17").Value = 10.89108910891
$ws.Range("M17").Value = 183.544303797468
$ws.Range("N17").Value = -16.728624535316

# Row 18
$ws.Range("C18").Value = 6
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 29
$ws.Range("G18").Value = 30
$ws.Range("H18").Value = -3.333333333333
$ws.Range("I18").Value = 175
$ws.Range("J18").Value = 158
$ws.Range("K18").Value = 10.759493670886
$ws.Range("L18").Value = -5.913978494623
$ws.Range("M18").Value = 31.578947368421
$ws.Range("N18").Value = -83.52165725047

# Row 19
$ws.Range("C19").Value = 35
$ws.Range("D19").Value = 38
$ws.Range("E19").Value = -7.894736842105
$ws.Range("F19").Value = 133
$ws.Range("G19").Value = 163
$ws.Range("H19").Value = -18.40490797546
$ws.Range("I19").Value = 702
$ws.Range("J19").Value = 841
$ws.Range("K19").Value = -16.527942925089
$ws.Range("L19").Value = -28.658536585365
$ws.Range("M19").Value = -22
$ws.Range("N19").Value = -81.408898305084

# Row 20
$ws.Range("F20").Value = 2
$ws.Range("H20").Value = -66.666666666666
$ws.Range("I20").Value = 10
$ws.Range("J20").Value = 21
$ws.Range("K20").Value = -52.380952380952
$ws.Range("L20").Value = -60
$ws.Range("M20").Value = 11.111111111111
$ws.Range("N20").Value = -93.670886075949

# Row 21
$ws.Range("C21").Value = 57
$ws.Range("D21").Value = 66
$ws.Range("E21").Value = -13.636363636363
$ws.Range("F21").Value = 240
$ws.Range("G21").Value = 276
$ws.Range("H21").Value = -13.043478260869
$ws.Range("I21").Value = 1269
$ws.Range("J21").Value = 1403
$ws.Range("K21").Value = -9.550962223806
$ws.Range("L21").Value = -20.934579439252
$ws.Range("M21").Value = 7.269653423499
$ws.Range("N21").Value = -79.812281259942

# Row 22
$ws.Range("C22").Value = 6
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = 200
$ws.Range("F22").Value = 19
$ws.Range("H22").Value = 18.75
$ws.Range("I22").Value = 92
$ws.Range("J22").Value = 75
$ws.Range("K22").Value = 22.666666666666
$ws.Range("L22").Value = -7.070707070707
$ws.Range("M22").Value = 80.392156862745

# Row 24
$ws.Range("D24").Value = 96
$ws.Range("E24").Value = -15.625
$ws.Range("F24").Value = 316
$ws.Range("G24").Value = 352
$ws.Range("H24").Value = -10.227272727272
$ws.Range("I24").Value = 1602
$ws.Range("J24").Value = 1692
$ws.Range("K24").Value = -5.31914893617
$ws.Range("L24").Value = 7.444668008048
$ws.Range("M24").Value = -12.459016393442

# Row 25
$ws.Range("C25").Value = 69
$ws.Range("D25").Value = 83
$ws.Range("E25").Value = -16.867469879518
$ws.Range("F25").Value = 268
$ws.Range("G25").Value = 304
$ws.Range("H25").Value = -11.842105263157
$ws.Range("I25").Value = 1401
$ws.Range("J25").Value = 1489
$ws.Range("K25").Value = -5.910006715916
$ws.Range("L25").Value = -0.284697508896

# Row 26
$ws.Range("C26").Value = 29
$ws.Range("D26").Value = 14
$ws.Range("E26").Value = 107.142857142857
$ws.Range("F26").Value = 90
$ws.Range("G26").Value = 73
$ws.Range("H26").Value = 23.287671232876
$ws.Range("I26").Value = 433
$ws.Range("J26").Value = 393
$ws.Range("K26").Value = 10.178117048346
$ws.Range("L26").Value = 10.178117048346
$ws.Range("M26").Value = 79.668049792531

# Row 27
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 2
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("E27").Value = 0
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 66.666666666666
$ws.Range("I27").Value = 21
$ws.Range("J27").Value = 18
$ws.Range("K27").Value = 16.666666666666
$ws.Range("L27").Value = 90.90909090909

# Row 28
$ws.Range("C28").Value = 7
$ws.Range("D28").Value = 4
$ws.Range("E28").Value = 75
$ws.Range("F28").Value = 33
$ws.Range("G28").Value = 12
$ws.Range("H28").Value = 175
$ws.Range("I28").Value = 96
$ws.Range("J28").Value = 81
$ws.Range("K28").Value = 18.518518518518
$ws.Range("L28").Value = 10.344827586206

# Row 31
$ws.Range("D31").Value = 1
$ws.Range("D31").NumberFormat = '#,##0'
$ws.Range("E31").Value = -100
$ws.Range("E31").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F31").Value = 2
$ws.Range("G31").Value = 2
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 9
$ws.Range("K31").Value = -11.111111111111
$ws.Range("L31").Value = 0

